$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "26.041.39"
$ws.Range('E2').Value = "  -2.11%  "
$ws.Range('D3').Value = "1.667.39"
$ws.Range('E3').Value = "  -1.60%  "
$ws.Range('D4').Value = "'1.004"
$ws.Range('E4').Value = "  -0.18%  "
$ws.Range('D5').Value = "'216.78"
$ws.Range('E5').Value = "  -1.37%  "
$ws.Range('D6').Value = "'0.5101"
$ws.Range('E6').Value = "  -0.24%  "
$ws.Range('E7').Value = "  -0.17%  "
$ws.Range('E8').Value = "  +0.11%  "
$ws.Range('D9').Value = "'0.06411"
$ws.Range('E9').Value = "  +1.37%  "
$ws.Range('D10').Value = "'21.84"
$ws.Range('E10').Value = "  -1.41%  "
$ws.Range('E11').Value = "  +1.23%  "
$ws.Range('B12').Value = "Polkadot"
$ws.Range('C12').Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('D12').Value = "'4.509"
$ws.Range('E12').Value = "  -0.46%  "
$ws.Range('B13').Value = "WrappedEther"
$ws.Range('C13').Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D13').Value = "1.668.46"
$ws.Range('E13').Value = "  -1.59%  "
$ws.Range('D14').Value = "'0.5838"
$ws.Range('E14').Value = "  +0.81%  "
$ws.Range('D15').Value = "'0.000008573"
$ws.Range('E15').Value = "  +0.66%  "
$ws.Range('D16').Value = "'64.33"
$ws.Range('E16').Value = "  -1.77%  "
$ws.Range('D17').Value = "26.126.05"
$ws.Range('E17').Value = "  -1.91%  "
$ws.Range('D18').Value = "'4.939"
$ws.Range('E18').Value = "  -0.85%  "
$ws.Range('D19').Value = "'1.004"
$ws.Range('E19').Value = "  -0.23%  "
$ws.Range('E20').Value = "  -1.72%  "
$ws.Range('D21').Value = "'191.95"
$ws.Range('E21').Value = "  +2.86%  "
$ws.Range('D22').Value = "'6.210"
$ws.Range('E22').Value = "  -0.73%  "
$ws.Range('E23').Value = "  -0.16%  "
$ws.Range('D24').Value = "'144.77"
$ws.Range('E24').Value = "  +0.06%  "
$ws.Range('D25').Value = "'7.621"
$ws.Range('E25').Value = "  +2.03%  "
$ws.Range('D26').Value = "'0.1197"
$ws.Range('E26').Value = "  +2.85%  "
$ws.Range('D27').Value = "'15.68"
$ws.Range('E27').Value = "  -1.02%  "
$ws.Range('D28').Value = "'0.06509"
$ws.Range('E28').Value = "  +13.04%  "
$ws.Range('D29').Value = "'1.326"
$ws.Range('E29').Value = "  -1.20%  "
$ws.Range('D30').Value = "'1.317"
$ws.Range('E30').Value = "  -1.91%  "
$ws.Range('D31').Value = "'3.541"
$ws.Range('E31').Value = "  +0.57%  "
$ws.Range('E32').Value = "  +0.34%  "
$ws.Range('E33').Value = "  +0.37%  "
$ws.Range('E34').Value = "  +0.04%  "
$ws.Range('D35').Value = "'0.6117"
$ws.Range('E35').Value = "  +2.01%  "
$ws.Range('D36').Value = "'2.370"
$ws.Range('E36').Value = "  +0.39%  "
$ws.Range('D37').Value = "'2.682"
$ws.Range('E37').Value = "  -0.04%  "
$ws.Range('D38').Value = "'6.274"
$ws.Range('E38').Value = "  +7.52%  "
$ws.Range('D39').Value = "'0.01602"
$ws.Range('E39').Value = "  -1.14%  "
$ws.Range('D40').Value = "1.090.34"
$ws.Range('E40').Value = "  -1.16%  "
$ws.Range('D41').Value = "'0.8636"
$ws.Range('E41').Value = "  +0.73%  "
$ws.Range('E42').Value = "  +0.57%  "
$ws.Range('D43').Value = "'100.90"
$ws.Range('E43').Value = "  +1.54%  "
$ws.Range('D44').Value = "1.816.81"
$ws.Range('E44').Value = "  -1.86%  "
$ws.Range('E45').Value = "  -2.24%  "
$ws.Range('E46').Value = "  -0.21%  "
$ws.Range('D47').Value = "'1.007"
$ws.Range('E47').Value = "  +0.10%  "
$ws.Range('D48').Value = "'8.049"
$ws.Range('E48').Value = "  -0.74%  "
$ws.Range('D49').Value = "'0.05230"
$ws.Range('E49').Value = "  -0.06%  "
$ws.Range('D50').Value = "'0.4284"
$ws.Range('E50').Value = "  -1.02%  "
$ws.Range('D51').Value = "'6.061"
$ws.Range('E51').Value = "  +4.53%  "
